$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three changed data cells on row 2
# D2 must stay a text value (not get converted to a number), so prefix it
# with an apostrophe to force Excel to keep it as text.
$ws.Range("D2").Value = "'7"
$ws.Range("G2").Value = "AAACT23201MY9ZV0Q"
$ws.Range("H2").Value = "20 jul. 2023, 09:31:00"

# Update the active selection to match the saved state
$ws.Range("H7").Select()
